# Weekly CompStat update (34th Precinct): new crime data collected.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Header block: new Mayor name, new volume/issue number, new report week.
# ---------------------------------------------------------------------------
$ws.Range("M6").Value = "Thomas G. Donlon"
$ws.Range("A8").Value = "Volume 31   Number  39"
$ws.Range("C9").Value = "Report Covering the Week  9/23/2024  Through  9/29/2024"

# ---------------------------------------------------------------------------
# Column H got a bit wider (bestFit recalculated after the data refresh).
# ---------------------------------------------------------------------------
$ws.Columns.Item(8).ColumnWidth = 6.7

# ---------------------------------------------------------------------------
# Helper: some cells flip between a numeric value and the text placeholder
# "0" (shared string) used when a figure is not meaningful, or vice-versa.
# Re-use an already-correctly-styled neighbour cell's format so the
# underlying cell style index is preserved/matched.
# ---------------------------------------------------------------------------
function Set-TextZero($target, $formatSource) {
    $ws.Range($target).Value = "'0"
    $ws.Range($formatSource).Copy()
    $ws.Range($target).PasteSpecial(-4122)
}
function Set-NumberLike($target, $value, $formatSource) {
    $ws.Range($target).Value = $value
    $ws.Range($formatSource).Copy()
    $ws.Range($target).PasteSpecial(-4122)
}

# Row 14 - Murder
Set-TextZero "C14" "D14"
$ws.Range("N14").Value = -87.755102040816

# Row 15 - Rape
$ws.Range("G15").Value = 1
$ws.Range("H15").Value = 300
$ws.Range("M15").Value = -7.692307692307
$ws.Range("N15").Value = -83.098591549295

# Row 16 - Robbery
$ws.Range("C16").Value = 3
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 13
$ws.Range("G16").Value = 16
$ws.Range("H16").Value = -18.75
$ws.Range("I16").Value = 174
$ws.Range("J16").Value = 163
$ws.Range("K16").Value = 6.748466257668
$ws.Range("L16").Value = -16.746411483253
$ws.Range("M16").Value = -13.432835820895
$ws.Range("N16").Value = -84.281842818428

# Row 17 - Fel. Assault
$ws.Range("C17").Value = 4
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = -20
$ws.Range("F17").Value = 17
$ws.Range("H17").Value = -32
$ws.Range("I17").Value = 256
$ws.Range("J17").Value = 248
$ws.Range("K17").Value = 3.225806451612
$ws.Range("L17").Value = 6.224066390041
$ws.Range("M17").Value = 56.097560975609
$ws.Range("N17").Value = -64.738292011019

# Row 18 - Burglary
$ws.Range("C18").Value = 4
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = 33.333333333333
$ws.Range("G18").Value = 12
$ws.Range("H18").Value = -16.666666666666
$ws.Range("I18").Value = 86
$ws.Range("J18").Value = 116
$ws.Range("K18").Value = -25.862068965517
$ws.Range("L18").Value = -23.893805309734
$ws.Range("M18").Value = -35.820895522388
$ws.Range("N18").Value = -94.707692307692

# Row 19 - Gr. Larceny
$ws.Range("C19").Value = 12
$ws.Range("D19").Value = 14
$ws.Range("E19").Value = -14.285714285714
$ws.Range("F19").Value = 54
$ws.Range("G19").Value = 38
$ws.Range("H19").Value = 42.105263157894
$ws.Range("I19").Value = 475
$ws.Range("J19").Value = 409
$ws.Range("K19").Value = 16.136919315403
$ws.Range("L19").Value = 8.447488584474
$ws.Range("M19").Value = 79.245283018867
$ws.Range("N19").Value = -51.874366767983

# Row 20 - G.L.A.
$ws.Range("C20").Value = 8
$ws.Range("D20").Value = 5
$ws.Range("E20").Value = 60
$ws.Range("F20").Value = 11
$ws.Range("G20").Value = 13
$ws.Range("H20").Value = -15.384615384615
$ws.Range("I20").Value = 129
$ws.Range("J20").Value = 173
$ws.Range("K20").Value = -25.43352601156
$ws.Range("L20").Value = -46.913580246913
$ws.Range("M20").Value = 98.461538461538
$ws.Range("N20").Value = -91.342281879194

# Row 21 - TOTAL
$ws.Range("C21").Value = 32
$ws.Range("D21").Value = 30
$ws.Range("E21").Value = 6.666666666666
$ws.Range("F21").Value = 110
$ws.Range("G21").Value = 105
$ws.Range("H21").Value = 4.761904761904
$ws.Range("I21").Value = 1138
$ws.Range("J21").Value = 1118
$ws.Range("K21").Value = 1.788908765652
$ws.Range("L21").Value = -10.181531176006
$ws.Range("M21").Value = 34.834123222748
$ws.Range("N21").Value = -81.205615194054

# Row 22 - Transit
Set-NumberLike "D22" 2 "G15"
Set-NumberLike "E22" -100 "H15"
Set-TextZero "F22" "D14"
$ws.Range("G22").Value = 3
$ws.Range("H22").Value = -100
$ws.Range("J22").Value = 31
$ws.Range("K22").Value = -38.709677419354

# Row 23 - Housing
Set-NumberLike "C23" 1 "D23"
$ws.Range("D23").Value = 1
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = 3
$ws.Range("H23").Value = -40
$ws.Range("I23").Value = 33
$ws.Range("J23").Value = 26
$ws.Range("K23").Value = 26.923076923076
$ws.Range("L23").Value = 10
$ws.Range("M23").Value = 120

# Row 24 - Petit Larceny
$ws.Range("C24").Value = 26
$ws.Range("D24").Value = 19
$ws.Range("E24").Value = 36.842105263157
$ws.Range("F24").Value = 124
$ws.Range("G24").Value = 118
$ws.Range("H24").Value = 5.084745762711
$ws.Range("I24").Value = 1104
$ws.Range("J24").Value = 1051
$ws.Range("K24").Value = 5.042816365366
$ws.Range("L24").Value = 2.033271719038
$ws.Range("M24").Value = 139.479392624729

# Row 25 - Retail Theft
$ws.Range("C25").Value = 13
$ws.Range("D25").Value = 3
$ws.Range("E25").Value = 333.333333333333
$ws.Range("F25").Value = 49
$ws.Range("G25").Value = 56
$ws.Range("H25").Value = -12.5
$ws.Range("I25").Value = 528
$ws.Range("J25").Value = 475
$ws.Range("K25").Value = 11.157894736842
$ws.Range("L25").Value = -4.174228675136

# Row 26 - Misd. Assault
$ws.Range("C26").Value = 14
$ws.Range("D26").Value = 10
$ws.Range("E26").Value = 40
$ws.Range("F26").Value = 38
$ws.Range("G26").Value = 49
$ws.Range("H26").Value = -22.448979591836
$ws.Range("I26").Value = 433
$ws.Range("J26").Value = 403
$ws.Range("K26").Value = 7.444168734491
$ws.Range("L26").Value = 26.239067055393
$ws.Range("M26").Value = 3.836930455635

# Row 27 - UCR Rape*
$ws.Range("G27").Value = 2
$ws.Range("H27").Value = 100

# Row 28 - Other Sex Crimes
$ws.Range("I28").Value = 45
$ws.Range("K28").Value = 4.651162790697
$ws.Range("L28").Value = 40.625

# Row 29 - Shooting Vic.
Set-TextZero "C29" "D29"
$ws.Range("L29").Value = 5.555555555555
$ws.Range("N29").Value = -87.662337662337

# Row 30 - Shooting Inc.
Set-TextZero "C30" "D30"
$ws.Range("L30").Value = -11.111111111111
$ws.Range("N30").Value = -88.235294117647
